$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header summary figures
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 981207      # VALOR MORA total
$ws.Range("F13").Value = 27          # Cant. Periodos

# ---------------------------------------------------------------------------
# 2) Rebuild the period rows (16..42) in descending (newest-first) order and
#    add the newest period (2507) that was not present before.
# ---------------------------------------------------------------------------

# First, give row 42 the same look (borders/shading) as the current last data
# row (41), and promote row 41 to a "normal" interior-row look (copied from
# row 40), since row 41 used to be the bottom border row and that role now
# belongs to the new row 42.
$ws.Range("B41:J41").Copy() | Out-Null
$ws.Range("B42:J42").PasteSpecial(-4122) | Out-Null
$ws.Range("B40:J40").Copy() | Out-Null
$ws.Range("B41:J41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the worker/base data for the newly-created row 42 (same worker as
# every other row).
$ws.Range("B42").Value = "CC"
$ws.Range("C42").Value = "1002185857"
$ws.Range("D42").Value = "JEAN PAUL CARDOZO DIAZ"
$ws.Range("F42").Value = 36341
$ws.Range("G42").Value = 908526

# Now write the period labels for every data row, newest period first.
$periods = @("2507","2506","2505","2504","2503","2502","2501", `
             "2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401", `
             "2312","2311","2310","2309","2308","2307","2306","2305")

$row = 16
foreach ($p in $periods) {
    $ws.Range("E" + $row).Value = $p
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3) Re-order the signature block: what used to be the underline row (46)
#    now sits at row 47, and the name/signature captions that used to be on
#    row 47 move down to the new row 48.
# ---------------------------------------------------------------------------
$ws.Rows("46").Insert()

$wb.Save()
